$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260, shifting existing rows 260:273 down to 261:274
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with the new weekly price record
$ws.Range("A260").Value = 5
$ws.Range("B260").Value = "Macroferia Regional de Talca"
$ws.Range("C260").Value = "Maule"
$ws.Range("D260").Value = 44753
$ws.Range("E260").Value = 7
$ws.Range("F260").Value = 100112009
$ws.Range("G260").Value = "Acelga"
$ws.Range("H260").Value = "Sin especificar"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 400
$ws.Range("K260").Value = 3000
$ws.Range("L260").Value = 3000
$ws.Range("M260").Value = 3000
$ws.Range("N260").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O260").Value = "Región del Maule"
$ws.Range("P260").Value = 750
$ws.Range("Q260").Value = 4
$ws.Range("R260").Value = "Hortaliza"
